$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grab a "header/label" style template cell before we start mutating ---
# A1 currently holds the old "section" header and carries the bold/bordered
# style (s="1") that is shared by the header row and the first data column.
$styleTemplate = $ws.Range("A1")
$styleTemplate.Copy()

# New header row (English). Column A no longer has a header (it becomes a
# plain numeric index), so headers now start at column B.
$headers = @("section_time","killed","severly_injured","lightly_injured","injured_pedestrians","casualties_ages_0-19","casualties_ages_20-64","casualties_ages_65_plus","total_casualties","vehicle_count","drivers")
for ($j = 0; $j -lt $headers.Length; $j++) {
    $ws.Cells.Item(1, $j + 2).Value = $headers[$j]
}

# Paste the bold/bordered style onto the new last header cell (L1) so it
# matches the rest of the header row.
$ws.Cells.Item(1, 12).PasteSpecial(-4122)

# New data rows: index, section_time, killed, severly_injured, lightly_injured,
# injured_pedestrians, casualties_ages_0-19, casualties_ages_20-64,
# casualties_ages_65_plus, total_casualties, vehicle_count, drivers
$data = @(
    @(0, "not_extended_after_2012", 24, 69, 202, 2, 53, 198, 35, 295, 135, 135),
    @(1, "not_extended_before_2012", 33, 62, 423, 0, 187, 314, 16, 518, 208, 208),
    @(2, "phase_1_after_2012", 2, 13, 91, 1, 27, 75, 2, 106, 60, 60),
    @(3, "phase_1_before_2012", 13, 15, 246, 1, 70, 197, 5, 274, 147, 147),
    @(4, "phase_2_after_2012", 7, 0, 28, 1, 10, 23, 2, 35, 18, 18),
    @(5, "phase_2_before_2012", 8, 17, 78, 0, 27, 74, 2, 103, 40, 40)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}

# Paste the bold/bordered style onto the new A7 cell (the extra row added by
# the split of "before/after 2012") so the first column keeps its style.
$ws.Cells.Item(7, 1).PasteSpecial(-4122)

# Finally, drop the old A1 "section" header entirely (no value, no style) -
# it was removed outright in favour of the new B1 "section_time" header.
$ws.Range("A1").Clear()

Write-Host "edit applied"
